$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 03:07"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8037789
$ws.Range("C4").Value = 45791
$ws.Range("D4").Value = 5184615
$ws.Range("E4").Value = 2633163
$ws.Range("G4").Value = 316
$ws.Range("H4").Value = 220011

# Row 5 - India
$ws.Range("B5").Value = 7173565
$ws.Range("C5").Value = 54265
$ws.Range("D5").Value = 6224792
$ws.Range("E5").Value = 838879

# Row 11 - Peru
$ws.Range("B11").Value = 851171
$ws.Range("C11").Value = 1800
$ws.Range("D11").Value = 748097
$ws.Range("E11").Value = 69717
$ws.Range("G11").Value = 52
$ws.Range("H11").Value = 33357

# Row 25 - Alemania
$ws.Range("B25").Value = 331094
$ws.Range("C25").Value = 4803
$ws.Range("E25").Value = 44473

# Row 155 - Guinea-Bisau
$ws.Range("B155").Value = 2389
$ws.Range("C155").Value = 4
$ws.Range("D155").Value = 1782
$ws.Range("E155").Value = 566
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 41

# Row 167 - Niger
$ws.Range("B167").Value = 1202
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 1124

# Row 169 - Santo Tome y Principe
$ws.Range("D169").Value = 895
$ws.Range("E169").Value = 19

# Row 192 - Bermudas
$ws.Range("D192").Value = 172
$ws.Range("E192").Value = 3

# Row 209 - Granada
$ws.Range("B209").Value = 25
$ws.Range("C209").Value = 1
$ws.Range("E209").Value = 1
